# New PO forecast model
# - Appends newly observed actuals to "Weekly Quantity" and "Monthly Trend"
# - Replaces the forecast series on "PO Forecast" with the output of the
#   refreshed model (updated historical values + new future periods)

$wb = $excel.ActiveWorkbook

function Set-DateRow {
    param($ws, $r, $dateVal, $qty)
    $ws.Cells.Item($r, 1).Value = $dateVal
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item($r - 1, 1).NumberFormat
    $ws.Cells.Item($r, 2).Value = $qty
}

# ---------------------------------------------------------------------------
# Sheet: Weekly Quantity  (append rows 21-23)
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

Set-DateRow $wsWeekly 21 45662.99999999999 48
Set-DateRow $wsWeekly 22 45669.99999999999 24
Set-DateRow $wsWeekly 23 45676.99999999999 4

# ---------------------------------------------------------------------------
# Sheet: Monthly Trend  (append row 12)
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

Set-DateRow $wsMonthly 12 45688.99999999999 76

# ---------------------------------------------------------------------------
# Sheet: PO Forecast  (refresh forecast values for rows 2-28, add rows 29-31)
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Rows whose date (column A) stays the same but forecast qty (column B) changes
$wsForecast.Cells.Item(2, 2).Value = 23
$wsForecast.Cells.Item(3, 2).Value = 30
$wsForecast.Cells.Item(4, 2).Value = 31
$wsForecast.Cells.Item(5, 2).Value = 32
$wsForecast.Cells.Item(6, 2).Value = 34
$wsForecast.Cells.Item(7, 2).Value = 35
$wsForecast.Cells.Item(8, 2).Value = 37
$wsForecast.Cells.Item(9, 2).Value = 38
$wsForecast.Cells.Item(10, 2).Value = 40
$wsForecast.Cells.Item(11, 2).Value = 41
$wsForecast.Cells.Item(12, 2).Value = 41
$wsForecast.Cells.Item(13, 2).Value = 42
$wsForecast.Cells.Item(14, 2).Value = 42
$wsForecast.Cells.Item(15, 2).Value = 43
$wsForecast.Cells.Item(16, 2).Value = 43
$wsForecast.Cells.Item(17, 2).Value = 44
$wsForecast.Cells.Item(18, 2).Value = 46
$wsForecast.Cells.Item(19, 2).Value = 46
$wsForecast.Cells.Item(20, 2).Value = 47

# Rows whose date (column A) AND forecast qty (column B) both change
$wsForecast.Cells.Item(21, 1).Value = 45662.99999999999
$wsForecast.Cells.Item(21, 2).Value = 48

$wsForecast.Cells.Item(22, 1).Value = 45669.99999999999
$wsForecast.Cells.Item(22, 2).Value = 49

$wsForecast.Cells.Item(23, 1).Value = 45676.99999999999
$wsForecast.Cells.Item(23, 2).Value = 50

$wsForecast.Cells.Item(24, 1).Value = 45683.99999999999
$wsForecast.Cells.Item(24, 2).Value = 50

$wsForecast.Cells.Item(25, 1).Value = 45690.99999999999
$wsForecast.Cells.Item(25, 2).Value = 51

$wsForecast.Cells.Item(26, 1).Value = 45697.99999999999
$wsForecast.Cells.Item(26, 2).Value = 51

$wsForecast.Cells.Item(27, 1).Value = 45704.99999999999
$wsForecast.Cells.Item(27, 2).Value = 52

$wsForecast.Cells.Item(28, 1).Value = 45711.99999999999
$wsForecast.Cells.Item(28, 2).Value = 52

# New trailing forecast rows
Set-DateRow $wsForecast 29 45718.99999999999 53
Set-DateRow $wsForecast 30 45725.99999999999 53
Set-DateRow $wsForecast 31 45732.99999999999 54
